$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add CT1/CU1 (copy header style s=1 from CR1, then set text) ---
$ws.Range("CR1").Copy()
$ws.Range("CT1").PasteSpecial(-4122)
$ws.Range("CU1").PasteSpecial(-4122)
$ws.Range("CT1").Value = "07-09_A"
$ws.Range("CU1").Value = "07-09_0"

# --- Fix A172 / A173 to numeric (were inline strings in the source) ---
$ws.Range("A172").Value = 59789149
$ws.Range("A173").Value = 59742804

# --- Data rows: CS -> numeric (same value), add CT (colored numeric) and CU (text) ---
$ws.Range("CS2").Value = 2494
$ws.Range("CT2").Value = 0
$ws.Range("CT2").Interior.Color = 255
$ws.Range("CU2").NumberFormat = "@"
$ws.Range("CU2").Value = "2494"
$ws.Range("CS3").Value = 0
$ws.Range("CT3").Value = 0
$ws.Range("CT3").Interior.Color = 255
$ws.Range("CU3").NumberFormat = "@"
$ws.Range("CU3").Value = "0"
$ws.Range("CS4").Value = 0
$ws.Range("CT4").Value = 0
$ws.Range("CT4").Interior.Color = 255
$ws.Range("CU4").NumberFormat = "@"
$ws.Range("CU4").Value = "0"
$ws.Range("CS5").Value = 0
$ws.Range("CT5").Value = 0
$ws.Range("CT5").Interior.Color = 255
$ws.Range("CU5").NumberFormat = "@"
$ws.Range("CU5").Value = "0"
$ws.Range("CS6").Value = 0
$ws.Range("CT6").Value = 0
$ws.Range("CT6").Interior.Color = 255
$ws.Range("CU6").NumberFormat = "@"
$ws.Range("CU6").Value = "0"
$ws.Range("CS18").Value = 3117
$ws.Range("CT18").Value = 20
$ws.Range("CT18").Interior.Color = 16777215
$ws.Range("CU18").NumberFormat = "@"
$ws.Range("CU18").Value = "3721"
$ws.Range("CS20").Value = 3902
$ws.Range("CT20").Value = 20
$ws.Range("CT20").Interior.Color = 16777215
$ws.Range("CU20").NumberFormat = "@"
$ws.Range("CU20").Value = "4138"
$ws.Range("CS21").Value = 2667
$ws.Range("CT21").Value = 0
$ws.Range("CT21").Interior.Color = 255
$ws.Range("CU21").NumberFormat = "@"
$ws.Range("CU21").Value = "2718"
$ws.Range("CS22").Value = 4273
$ws.Range("CT22").Value = 21
$ws.Range("CT22").Interior.Color = 16777215
$ws.Range("CU22").NumberFormat = "@"
$ws.Range("CU22").Value = "4713"
$ws.Range("CS23").Value = 4660
$ws.Range("CT23").Value = 23
$ws.Range("CT23").Interior.Color = 16777215
$ws.Range("CU23").NumberFormat = "@"
$ws.Range("CU23").Value = "5159"
$ws.Range("CS24").Value = 4471
$ws.Range("CT24").Value = 35
$ws.Range("CT24").Interior.Color = 32768
$ws.Range("CU24").NumberFormat = "@"
$ws.Range("CU24").Value = "4833"
$ws.Range("CS27").Value = 0
$ws.Range("CT27").Value = 0
$ws.Range("CT27").Interior.Color = 255
$ws.Range("CU27").NumberFormat = "@"
$ws.Range("CU27").Value = "0"
$ws.Range("CS29").Value = 2820
$ws.Range("CT29").Value = 0
$ws.Range("CT29").Interior.Color = 255
$ws.Range("CU29").NumberFormat = "@"
$ws.Range("CU29").Value = "2836"
$ws.Range("CS30").Value = 4270
$ws.Range("CT30").Value = 0
$ws.Range("CT30").Interior.Color = 255
$ws.Range("CU30").NumberFormat = "@"
$ws.Range("CU30").Value = "4275"
$ws.Range("CS31").Value = 4421
$ws.Range("CT31").Value = 30
$ws.Range("CT31").Interior.Color = 16777215
$ws.Range("CU31").NumberFormat = "@"
$ws.Range("CU31").Value = "4837"
$ws.Range("CS32").Value = 2577
$ws.Range("CT32").Value = 0
$ws.Range("CT32").Interior.Color = 255
$ws.Range("CU32").NumberFormat = "@"
$ws.Range("CU32").Value = "2608"
$ws.Range("CS34").Value = 0
$ws.Range("CT34").Value = 0
$ws.Range("CT34").Interior.Color = 255
$ws.Range("CU34").NumberFormat = "@"
$ws.Range("CU34").Value = "0"
$ws.Range("CS35").Value = 0
$ws.Range("CT35").Value = 0
$ws.Range("CT35").Interior.Color = 255
$ws.Range("CU35").NumberFormat = "@"
$ws.Range("CU35").Value = "0"
$ws.Range("CS36").Value = 2635
$ws.Range("CT36").Value = 0
$ws.Range("CT36").Interior.Color = 255
$ws.Range("CU36").NumberFormat = "@"
$ws.Range("CU36").Value = "2702"
$ws.Range("CS38").Value = 4108
$ws.Range("CT38").Value = 22
$ws.Range("CT38").Interior.Color = 16777215
$ws.Range("CU38").NumberFormat = "@"
$ws.Range("CU38").Value = "4725"
$ws.Range("CS39").Value = 4116
$ws.Range("CT39").Value = 29
$ws.Range("CT39").Interior.Color = 16777215
$ws.Range("CU39").NumberFormat = "@"
$ws.Range("CU39").Value = "4449"
$ws.Range("CS41").Value = 2864
$ws.Range("CT41").Value = 20
$ws.Range("CT41").Interior.Color = 16777215
$ws.Range("CU41").NumberFormat = "@"
$ws.Range("CU41").Value = "3605"
$ws.Range("CS42").Value = 2601
$ws.Range("CT42").Value = 0
$ws.Range("CT42").Interior.Color = 255
$ws.Range("CU42").NumberFormat = "@"
$ws.Range("CU42").Value = "2664"
$ws.Range("CS46").Value = 4100
$ws.Range("CT46").Value = 20
$ws.Range("CT46").Interior.Color = 16777215
$ws.Range("CU46").NumberFormat = "@"
$ws.Range("CU46").Value = "4382"
$ws.Range("CS47").Value = 4827
$ws.Range("CT47").Value = 33
$ws.Range("CT47").Interior.Color = 32768
$ws.Range("CU47").NumberFormat = "@"
$ws.Range("CU47").Value = "5142"
$ws.Range("CS48").Value = 0
$ws.Range("CT48").Value = 0
$ws.Range("CT48").Interior.Color = 255
$ws.Range("CU48").NumberFormat = "@"
$ws.Range("CU48").Value = "0"
$ws.Range("CS49").Value = 4197
$ws.Range("CT49").Value = 27
$ws.Range("CT49").Interior.Color = 16777215
$ws.Range("CU49").NumberFormat = "@"
$ws.Range("CU49").Value = "4520"
$ws.Range("CS50").Value = 4296
$ws.Range("CT50").Value = 23
$ws.Range("CT50").Interior.Color = 16777215
$ws.Range("CU50").NumberFormat = "@"
$ws.Range("CU50").Value = "4704"
$ws.Range("CS52").Value = 4563
$ws.Range("CT52").Value = 30
$ws.Range("CT52").Interior.Color = 16777215
$ws.Range("CU52").NumberFormat = "@"
$ws.Range("CU52").Value = "4867"
$ws.Range("CS53").Value = 3292
$ws.Range("CT53").Value = 5
$ws.Range("CT53").Interior.Color = 65535
$ws.Range("CU53").NumberFormat = "@"
$ws.Range("CU53").Value = "3460"
$ws.Range("CS55").Value = 3013
$ws.Range("CT55").Value = 20
$ws.Range("CT55").Interior.Color = 16777215
$ws.Range("CU55").NumberFormat = "@"
$ws.Range("CU55").Value = "3546"
$ws.Range("CS56").Value = 4665
$ws.Range("CT56").Value = 30
$ws.Range("CT56").Interior.Color = 16777215
$ws.Range("CU56").NumberFormat = "@"
$ws.Range("CU56").Value = "4874"
$ws.Range("CS57").Value = 4058
$ws.Range("CT57").Value = 20
$ws.Range("CT57").Interior.Color = 16777215
$ws.Range("CU57").NumberFormat = "@"
$ws.Range("CU57").Value = "4174"
$ws.Range("CS58").Value = 3995
$ws.Range("CT58").Value = 20
$ws.Range("CT58").Interior.Color = 16777215
$ws.Range("CU58").NumberFormat = "@"
$ws.Range("CU58").Value = "4193"
$ws.Range("CS59").Value = 3900
$ws.Range("CT59").Value = 20
$ws.Range("CT59").Interior.Color = 16777215
$ws.Range("CU59").NumberFormat = "@"
$ws.Range("CU59").Value = "4008"
$ws.Range("CS60").Value = 4236
$ws.Range("CT60").Value = 27
$ws.Range("CT60").Interior.Color = 16777215
$ws.Range("CU60").NumberFormat = "@"
$ws.Range("CU60").Value = "4413"
$ws.Range("CS62").Value = 3851
$ws.Range("CT62").Value = 30
$ws.Range("CT62").Interior.Color = 16777215
$ws.Range("CU62").NumberFormat = "@"
$ws.Range("CU62").Value = "3989"
$ws.Range("CS63").Value = 3806
$ws.Range("CT63").Value = 6
$ws.Range("CT63").Interior.Color = 65535
$ws.Range("CU63").NumberFormat = "@"
$ws.Range("CU63").Value = "4002"
$ws.Range("CS64").Value = 3099
$ws.Range("CT64").Value = 0
$ws.Range("CT64").Interior.Color = 255
$ws.Range("CU64").NumberFormat = "@"
$ws.Range("CU64").Value = "3308"
$ws.Range("CS66").Value = 0
$ws.Range("CT66").Value = 0
$ws.Range("CT66").Interior.Color = 255
$ws.Range("CU66").NumberFormat = "@"
$ws.Range("CU66").Value = "0"
$ws.Range("CS67").Value = 0
$ws.Range("CT67").Value = 0
$ws.Range("CT67").Interior.Color = 255
$ws.Range("CU67").NumberFormat = "@"
$ws.Range("CU67").Value = "0"
$ws.Range("CS68").Value = 0
$ws.Range("CT68").Value = 0
$ws.Range("CT68").Interior.Color = 255
$ws.Range("CU68").NumberFormat = "@"
$ws.Range("CU68").Value = "0"
$ws.Range("CS70").Value = 0
$ws.Range("CT70").Value = 0
$ws.Range("CT70").Interior.Color = 255
$ws.Range("CU70").NumberFormat = "@"
$ws.Range("CU70").Value = "0"
$ws.Range("CS71").Value = 0
$ws.Range("CT71").Value = 0
$ws.Range("CT71").Interior.Color = 255
$ws.Range("CU71").NumberFormat = "@"
$ws.Range("CU71").Value = "0"
$ws.Range("CS72").Value = 0
$ws.Range("CT72").Value = 0
$ws.Range("CT72").Interior.Color = 255
$ws.Range("CU72").NumberFormat = "@"
$ws.Range("CU72").Value = "0"
$ws.Range("CS73").Value = 0
$ws.Range("CT73").Value = 0
$ws.Range("CT73").Interior.Color = 255
$ws.Range("CU73").NumberFormat = "@"
$ws.Range("CU73").Value = "0"
$ws.Range("CS74").Value = 0
$ws.Range("CT74").Value = 0
$ws.Range("CT74").Interior.Color = 255
$ws.Range("CU74").NumberFormat = "@"
$ws.Range("CU74").Value = "0"
$ws.Range("CS75").Value = 0
$ws.Range("CT75").Value = 0
$ws.Range("CT75").Interior.Color = 255
$ws.Range("CU75").NumberFormat = "@"
$ws.Range("CU75").Value = "0"
$ws.Range("CS76").Value = 2731
$ws.Range("CT76").Value = 0
$ws.Range("CT76").Interior.Color = 255
$ws.Range("CU76").NumberFormat = "@"
$ws.Range("CU76").Value = "2802"
$ws.Range("CS77").Value = 2789
$ws.Range("CT77").Value = 0
$ws.Range("CT77").Interior.Color = 255
$ws.Range("CU77").NumberFormat = "@"
$ws.Range("CU77").Value = "3027"
$ws.Range("CS78").Value = 2527
$ws.Range("CT78").Value = 0
$ws.Range("CT78").Interior.Color = 255
$ws.Range("CU78").NumberFormat = "@"
$ws.Range("CU78").Value = "2647"
$ws.Range("CS79").Value = 0
$ws.Range("CT79").Value = 0
$ws.Range("CT79").Interior.Color = 255
$ws.Range("CU79").NumberFormat = "@"
$ws.Range("CU79").Value = "0"
$ws.Range("CS80").Value = 0
$ws.Range("CT80").Value = 0
$ws.Range("CT80").Interior.Color = 255
$ws.Range("CU80").NumberFormat = "@"
$ws.Range("CU80").Value = "0"
$ws.Range("CS81").Value = 0
$ws.Range("CT81").Value = 0
$ws.Range("CT81").Interior.Color = 255
$ws.Range("CU81").NumberFormat = "@"
$ws.Range("CU81").Value = "0"
$ws.Range("CS82").Value = 0
$ws.Range("CT82").Value = 0
$ws.Range("CT82").Interior.Color = 255
$ws.Range("CU82").NumberFormat = "@"
$ws.Range("CU82").Value = "0"
$ws.Range("CS83").Value = 0
$ws.Range("CT83").Value = 0
$ws.Range("CT83").Interior.Color = 255
$ws.Range("CU83").NumberFormat = "@"
$ws.Range("CU83").Value = "0"
$ws.Range("CS84").Value = 0
$ws.Range("CT84").Value = 0
$ws.Range("CT84").Interior.Color = 255
$ws.Range("CU84").NumberFormat = "@"
$ws.Range("CU84").Value = "0"
$ws.Range("CS85").Value = 0
$ws.Range("CT85").Value = 0
$ws.Range("CT85").Interior.Color = 255
$ws.Range("CU85").NumberFormat = "@"
$ws.Range("CU85").Value = "0"
$ws.Range("CS86").Value = 0
$ws.Range("CT86").Value = 0
$ws.Range("CT86").Interior.Color = 255
$ws.Range("CU86").NumberFormat = "@"
$ws.Range("CU86").Value = "0"
$ws.Range("CS87").Value = 0
$ws.Range("CT87").Value = 0
$ws.Range("CT87").Interior.Color = 255
$ws.Range("CU87").NumberFormat = "@"
$ws.Range("CU87").Value = "0"
$ws.Range("CS88").Value = 0
$ws.Range("CT88").Value = 0
$ws.Range("CT88").Interior.Color = 255
$ws.Range("CU88").NumberFormat = "@"
$ws.Range("CU88").Value = "0"
$ws.Range("CS89").Value = 0
$ws.Range("CT89").Value = 0
$ws.Range("CT89").Interior.Color = 255
$ws.Range("CU89").NumberFormat = "@"
$ws.Range("CU89").Value = "0"
$ws.Range("CS90").Value = 0
$ws.Range("CT90").Value = 0
$ws.Range("CT90").Interior.Color = 255
$ws.Range("CU90").NumberFormat = "@"
$ws.Range("CU90").Value = "0"
$ws.Range("CS91").Value = 0
$ws.Range("CT91").Value = 0
$ws.Range("CT91").Interior.Color = 255
$ws.Range("CU91").NumberFormat = "@"
$ws.Range("CU91").Value = "0"
$ws.Range("CS92").Value = 0
$ws.Range("CT92").Value = 0
$ws.Range("CT92").Interior.Color = 255
$ws.Range("CU92").NumberFormat = "@"
$ws.Range("CU92").Value = "0"
$ws.Range("CS93").Value = 0
$ws.Range("CT93").Value = 0
$ws.Range("CT93").Interior.Color = 255
$ws.Range("CU93").NumberFormat = "@"
$ws.Range("CU93").Value = "0"
$ws.Range("CS94").Value = 0
$ws.Range("CT94").Value = 0
$ws.Range("CT94").Interior.Color = 255
$ws.Range("CU94").NumberFormat = "@"
$ws.Range("CU94").Value = "0"
$ws.Range("CS95").Value = 0
$ws.Range("CT95").Value = 0
$ws.Range("CT95").Interior.Color = 255
$ws.Range("CU95").NumberFormat = "@"
$ws.Range("CU95").Value = "0"
$ws.Range("CS97").Value = 0
$ws.Range("CT97").Value = 0
$ws.Range("CT97").Interior.Color = 255
$ws.Range("CU97").NumberFormat = "@"
$ws.Range("CU97").Value = "0"
$ws.Range("CS98").Value = 0
$ws.Range("CT98").Value = 0
$ws.Range("CT98").Interior.Color = 255
$ws.Range("CU98").NumberFormat = "@"
$ws.Range("CU98").Value = "0"
$ws.Range("CS99").Value = 0
$ws.Range("CT99").Value = 0
$ws.Range("CT99").Interior.Color = 255
$ws.Range("CU99").NumberFormat = "@"
$ws.Range("CU99").Value = "0"
$ws.Range("CS100").Value = 0
$ws.Range("CT100").Value = 0
$ws.Range("CT100").Interior.Color = 255
$ws.Range("CU100").NumberFormat = "@"
$ws.Range("CU100").Value = "0"
$ws.Range("CS101").Value = 0
$ws.Range("CT101").Value = 0
$ws.Range("CT101").Interior.Color = 255
$ws.Range("CU101").NumberFormat = "@"
$ws.Range("CU101").Value = "0"
$ws.Range("CS102").Value = 0
$ws.Range("CT102").Value = 0
$ws.Range("CT102").Interior.Color = 255
$ws.Range("CU102").NumberFormat = "@"
$ws.Range("CU102").Value = "0"
$ws.Range("CS103").Value = 0
$ws.Range("CT103").Value = 0
$ws.Range("CT103").Interior.Color = 255
$ws.Range("CU103").NumberFormat = "@"
$ws.Range("CU103").Value = "0"
$ws.Range("CS104").Value = 0
$ws.Range("CT104").Value = 0
$ws.Range("CT104").Interior.Color = 255
$ws.Range("CU104").NumberFormat = "@"
$ws.Range("CU104").Value = "0"
$ws.Range("CS105").Value = 0
$ws.Range("CT105").Value = 0
$ws.Range("CT105").Interior.Color = 255
$ws.Range("CU105").NumberFormat = "@"
$ws.Range("CU105").Value = "0"
$ws.Range("CS115").Value = 4741
$ws.Range("CT115").Value = 30
$ws.Range("CT115").Interior.Color = 16777215
$ws.Range("CU115").NumberFormat = "@"
$ws.Range("CU115").Value = "5196"
$ws.Range("CS116").Value = 0
$ws.Range("CT116").Value = 0
$ws.Range("CT116").Interior.Color = 255
$ws.Range("CU116").NumberFormat = "@"
$ws.Range("CU116").Value = "0"
$ws.Range("CS118").Value = 3430
$ws.Range("CT118").Value = 20
$ws.Range("CT118").Interior.Color = 16777215
$ws.Range("CU118").NumberFormat = "@"
$ws.Range("CU118").Value = "4089"
$ws.Range("CS119").Value = 0
$ws.Range("CT119").Value = 0
$ws.Range("CT119").Interior.Color = 255
$ws.Range("CU119").NumberFormat = "@"
$ws.Range("CU119").Value = "0"
$ws.Range("CS120").Value = 0
$ws.Range("CT120").Value = 0
$ws.Range("CT120").Interior.Color = 255
$ws.Range("CU120").NumberFormat = "@"
$ws.Range("CU120").Value = "0"
$ws.Range("CS121").Value = 0
$ws.Range("CT121").Value = 0
$ws.Range("CT121").Interior.Color = 255
$ws.Range("CU121").NumberFormat = "@"
$ws.Range("CU121").Value = "0"
$ws.Range("CS122").Value = 0
$ws.Range("CT122").Value = 0
$ws.Range("CT122").Interior.Color = 255
$ws.Range("CU122").NumberFormat = "@"
$ws.Range("CU122").Value = "0"
$ws.Range("CS123").Value = 2610
$ws.Range("CT123").Value = 7
$ws.Range("CT123").Interior.Color = 65535
$ws.Range("CU123").NumberFormat = "@"
$ws.Range("CU123").Value = "2719"
$ws.Range("CS125").Value = 0
$ws.Range("CT125").Value = 0
$ws.Range("CT125").Interior.Color = 255
$ws.Range("CU125").NumberFormat = "@"
$ws.Range("CU125").Value = "0"
$ws.Range("CS126").Value = 0
$ws.Range("CT126").Value = 0
$ws.Range("CT126").Interior.Color = 255
$ws.Range("CU126").NumberFormat = "@"
$ws.Range("CU126").Value = "0"
$ws.Range("CS129").Value = 0
$ws.Range("CT129").Value = 0
$ws.Range("CT129").Interior.Color = 255
$ws.Range("CU129").NumberFormat = "@"
$ws.Range("CU129").Value = "0"
$ws.Range("CS131").Value = 2537
$ws.Range("CT131").Value = 0
$ws.Range("CT131").Interior.Color = 255
$ws.Range("CU131").NumberFormat = "@"
$ws.Range("CU131").Value = "2553"
$ws.Range("CS132").Value = 2866
$ws.Range("CT132").Value = 0
$ws.Range("CT132").Interior.Color = 255
$ws.Range("CU132").NumberFormat = "@"
$ws.Range("CU132").Value = "2936"
$ws.Range("CS133").Value = 2413
$ws.Range("CT133").Value = 0
$ws.Range("CT133").Interior.Color = 255
$ws.Range("CU133").NumberFormat = "@"
$ws.Range("CU133").Value = "2399"
$ws.Range("CS136").Value = 5336
$ws.Range("CT136").Value = 39
$ws.Range("CT136").Interior.Color = 32768
$ws.Range("CU136").NumberFormat = "@"
$ws.Range("CU136").Value = "5874"
$ws.Range("CS137").Value = 4810
$ws.Range("CT137").Value = 36
$ws.Range("CT137").Interior.Color = 32768
$ws.Range("CU137").NumberFormat = "@"
$ws.Range("CU137").Value = "5367"
$ws.Range("CS139").Value = 5042
$ws.Range("CT139").Value = 39
$ws.Range("CT139").Interior.Color = 32768
$ws.Range("CU139").NumberFormat = "@"
$ws.Range("CU139").Value = "5583"
$ws.Range("CS140").Value = 0
$ws.Range("CT140").Value = 0
$ws.Range("CT140").Interior.Color = 255
$ws.Range("CU140").NumberFormat = "@"
$ws.Range("CU140").Value = "0"
$ws.Range("CS142").Value = 2495
$ws.Range("CT142").Value = 0
$ws.Range("CT142").Interior.Color = 255
$ws.Range("CU142").NumberFormat = "@"
$ws.Range("CU142").Value = "2495"
$ws.Range("CS143").Value = 0
$ws.Range("CT143").Value = 0
$ws.Range("CT143").Interior.Color = 255
$ws.Range("CU143").NumberFormat = "@"
$ws.Range("CU143").Value = "0"
$ws.Range("CS144").Value = 0
$ws.Range("CT144").Value = 0
$ws.Range("CT144").Interior.Color = 255
$ws.Range("CU144").NumberFormat = "@"
$ws.Range("CU144").Value = "1500"
$ws.Range("CS145").Value = 0
$ws.Range("CT145").Value = 0
$ws.Range("CT145").Interior.Color = 255
$ws.Range("CU145").NumberFormat = "@"
$ws.Range("CU145").Value = "0"
$ws.Range("CS146").Value = 0
$ws.Range("CT146").Value = 0
$ws.Range("CT146").Interior.Color = 255
$ws.Range("CU146").NumberFormat = "@"
$ws.Range("CU146").Value = "0"
$ws.Range("CS147").Value = 4174
$ws.Range("CT147").Value = 24
$ws.Range("CT147").Interior.Color = 16777215
$ws.Range("CU147").NumberFormat = "@"
$ws.Range("CU147").Value = "4468"
$ws.Range("CS148").Value = 0
$ws.Range("CT148").Value = 0
$ws.Range("CT148").Interior.Color = 255
$ws.Range("CU148").NumberFormat = "@"
$ws.Range("CU148").Value = "0"
$ws.Range("CS150").Value = 2668
$ws.Range("CT150").Value = 0
$ws.Range("CT150").Interior.Color = 255
$ws.Range("CU150").NumberFormat = "@"
$ws.Range("CU150").Value = "2663"
$ws.Range("CS151").Value = 0
$ws.Range("CT151").Value = 0
$ws.Range("CT151").Interior.Color = 255
$ws.Range("CU151").NumberFormat = "@"
$ws.Range("CU151").Value = "0"
$ws.Range("CS152").Value = 0
$ws.Range("CT152").Value = 0
$ws.Range("CT152").Interior.Color = 255
$ws.Range("CU152").NumberFormat = "@"
$ws.Range("CU152").Value = "0"
$ws.Range("CS153").Value = 2488
$ws.Range("CT153").Value = 0
$ws.Range("CT153").Interior.Color = 255
$ws.Range("CU153").NumberFormat = "@"
$ws.Range("CU153").Value = "2475"
$ws.Range("CS154").Value = 0
$ws.Range("CT154").Value = 0
$ws.Range("CT154").Interior.Color = 255
$ws.Range("CU154").NumberFormat = "@"
$ws.Range("CU154").Value = "0"
$ws.Range("CS157").Value = 0
$ws.Range("CT157").Value = 0
$ws.Range("CT157").Interior.Color = 255
$ws.Range("CU157").NumberFormat = "@"
$ws.Range("CU157").Value = "0"
$ws.Range("CS158").Value = 0
$ws.Range("CT158").Value = 0
$ws.Range("CT158").Interior.Color = 255
$ws.Range("CU158").NumberFormat = "@"
$ws.Range("CU158").Value = "0"
$ws.Range("CS159").Value = 0
$ws.Range("CT159").Value = 0
$ws.Range("CT159").Interior.Color = 255
$ws.Range("CU159").NumberFormat = "@"
$ws.Range("CU159").Value = "0"
$ws.Range("CS160").Value = 3124
$ws.Range("CT160").Value = 20
$ws.Range("CT160").Interior.Color = 16777215
$ws.Range("CU160").NumberFormat = "@"
$ws.Range("CU160").Value = "3512"
$ws.Range("CS161").Value = 2557
$ws.Range("CT161").Value = 0
$ws.Range("CT161").Interior.Color = 255
$ws.Range("CU161").NumberFormat = "@"
$ws.Range("CU161").Value = "2625"
$ws.Range("CS162").Value = 2114
$ws.Range("CT162").Value = 3
$ws.Range("CT162").Interior.Color = 65535
$ws.Range("CU162").NumberFormat = "@"
$ws.Range("CU162").Value = "2164"
$ws.Range("CS163").Value = 2525
$ws.Range("CT163").Value = 0
$ws.Range("CT163").Interior.Color = 255
$ws.Range("CU163").NumberFormat = "@"
$ws.Range("CU163").Value = "2565"
$ws.Range("CS165").Value = 0
$ws.Range("CT165").Value = 0
$ws.Range("CT165").Interior.Color = 255
$ws.Range("CU165").NumberFormat = "@"
$ws.Range("CU165").Value = "0"
$ws.Range("CS167").Value = 1500
$ws.Range("CT167").Value = 0
$ws.Range("CT167").Interior.Color = 255
$ws.Range("CU167").NumberFormat = "@"
$ws.Range("CU167").Value = "1500"
$ws.Range("CS168").Value = 1429
$ws.Range("CT168").Value = 0
$ws.Range("CT168").Interior.Color = 255
$ws.Range("CU168").NumberFormat = "@"
$ws.Range("CU168").Value = "0"
$ws.Range("CS169").Value = 0
$ws.Range("CT169").Value = 0
$ws.Range("CT169").Interior.Color = 255
$ws.Range("CU169").NumberFormat = "@"
$ws.Range("CU169").Value = "0"
$ws.Range("CS170").Value = 0
$ws.Range("CT170").Value = 0
$ws.Range("CT170").Interior.Color = 255
$ws.Range("CU170").NumberFormat = "@"
$ws.Range("CU170").Value = "0"
$ws.Range("CS171").Value = 3374
$ws.Range("CT171").Value = 8
$ws.Range("CT171").Interior.Color = 65535
$ws.Range("CU171").NumberFormat = "@"
$ws.Range("CU171").Value = "3622"
$ws.Range("CS172").Value = 1373
$ws.Range("CT172").Value = 0
$ws.Range("CT172").Interior.Color = 255
$ws.Range("CU172").NumberFormat = "@"
$ws.Range("CU172").Value = "1393"
$ws.Range("CS173").Value = 1631
$ws.Range("CT173").Value = 0
$ws.Range("CT173").Interior.Color = 255
$ws.Range("CU173").NumberFormat = "@"
$ws.Range("CU173").Value = "1619"

# --- Blank rows: just extend with empty CT (white fill, style 4) / CU (default style) ---
$ws.Range("CT7").Interior.Color = 16777215
$ws.Range("CU7").Font.Bold = $false
$ws.Range("CT8").Interior.Color = 16777215
$ws.Range("CU8").Font.Bold = $false
$ws.Range("CT9").Interior.Color = 16777215
$ws.Range("CU9").Font.Bold = $false
$ws.Range("CT10").Interior.Color = 16777215
$ws.Range("CU10").Font.Bold = $false
$ws.Range("CT11").Interior.Color = 16777215
$ws.Range("CU11").Font.Bold = $false
$ws.Range("CT12").Interior.Color = 16777215
$ws.Range("CU12").Font.Bold = $false
$ws.Range("CT13").Interior.Color = 16777215
$ws.Range("CU13").Font.Bold = $false
$ws.Range("CT14").Interior.Color = 16777215
$ws.Range("CU14").Font.Bold = $false
$ws.Range("CT15").Interior.Color = 16777215
$ws.Range("CU15").Font.Bold = $false
$ws.Range("CT16").Interior.Color = 16777215
$ws.Range("CU16").Font.Bold = $false
$ws.Range("CT17").Interior.Color = 16777215
$ws.Range("CU17").Font.Bold = $false
$ws.Range("CT19").Interior.Color = 16777215
$ws.Range("CU19").Font.Bold = $false
$ws.Range("CT25").Interior.Color = 16777215
$ws.Range("CU25").Font.Bold = $false
$ws.Range("CT26").Interior.Color = 16777215
$ws.Range("CU26").Font.Bold = $false
$ws.Range("CT28").Interior.Color = 16777215
$ws.Range("CU28").Font.Bold = $false
$ws.Range("CT33").Interior.Color = 16777215
$ws.Range("CU33").Font.Bold = $false
$ws.Range("CT37").Interior.Color = 16777215
$ws.Range("CU37").Font.Bold = $false
$ws.Range("CT40").Interior.Color = 16777215
$ws.Range("CU40").Font.Bold = $false
$ws.Range("CT43").Interior.Color = 16777215
$ws.Range("CU43").Font.Bold = $false
$ws.Range("CT44").Interior.Color = 16777215
$ws.Range("CU44").Font.Bold = $false
$ws.Range("CT45").Interior.Color = 16777215
$ws.Range("CU45").Font.Bold = $false
$ws.Range("CT51").Interior.Color = 16777215
$ws.Range("CU51").Font.Bold = $false
$ws.Range("CT54").Interior.Color = 16777215
$ws.Range("CU54").Font.Bold = $false
$ws.Range("CT61").Interior.Color = 16777215
$ws.Range("CU61").Font.Bold = $false
$ws.Range("CT65").Interior.Color = 16777215
$ws.Range("CU65").Font.Bold = $false
$ws.Range("CT69").Interior.Color = 16777215
$ws.Range("CU69").Font.Bold = $false
$ws.Range("CT96").Interior.Color = 16777215
$ws.Range("CU96").Font.Bold = $false
$ws.Range("CT106").Interior.Color = 16777215
$ws.Range("CU106").Font.Bold = $false
$ws.Range("CT107").Interior.Color = 16777215
$ws.Range("CU107").Font.Bold = $false
$ws.Range("CT108").Interior.Color = 16777215
$ws.Range("CU108").Font.Bold = $false
$ws.Range("CT109").Interior.Color = 16777215
$ws.Range("CU109").Font.Bold = $false
$ws.Range("CT110").Interior.Color = 16777215
$ws.Range("CU110").Font.Bold = $false
$ws.Range("CT111").Interior.Color = 16777215
$ws.Range("CU111").Font.Bold = $false
$ws.Range("CT112").Interior.Color = 16777215
$ws.Range("CU112").Font.Bold = $false
$ws.Range("CT113").Interior.Color = 16777215
$ws.Range("CU113").Font.Bold = $false
$ws.Range("CT114").Interior.Color = 16777215
$ws.Range("CU114").Font.Bold = $false
$ws.Range("CT117").Interior.Color = 16777215
$ws.Range("CU117").Font.Bold = $false
$ws.Range("CT124").Interior.Color = 16777215
$ws.Range("CU124").Font.Bold = $false
$ws.Range("CT127").Interior.Color = 16777215
$ws.Range("CU127").Font.Bold = $false
$ws.Range("CT128").Interior.Color = 16777215
$ws.Range("CU128").Font.Bold = $false
$ws.Range("CT130").Interior.Color = 16777215
$ws.Range("CU130").Font.Bold = $false
$ws.Range("CT134").Interior.Color = 16777215
$ws.Range("CU134").Font.Bold = $false
$ws.Range("CT135").Interior.Color = 16777215
$ws.Range("CU135").Font.Bold = $false
$ws.Range("CT138").Interior.Color = 16777215
$ws.Range("CU138").Font.Bold = $false
$ws.Range("CT141").Interior.Color = 16777215
$ws.Range("CU141").Font.Bold = $false
$ws.Range("CT149").Interior.Color = 16777215
$ws.Range("CU149").Font.Bold = $false
$ws.Range("CT155").Interior.Color = 16777215
$ws.Range("CU155").Font.Bold = $false
$ws.Range("CT156").Interior.Color = 16777215
$ws.Range("CU156").Font.Bold = $false
$ws.Range("CT164").Interior.Color = 16777215
$ws.Range("CU164").Font.Bold = $false
$ws.Range("CT166").Interior.Color = 16777215
$ws.Range("CU166").Font.Bold = $false

# --- New row 174: brand-new player, first appears on 07-09 ---
$ws.Range("A7:CS7").Copy()
$ws.Range("A174:CS174").PasteSpecial(-4122)
$ws.Range("A174").NumberFormat = "@"
$ws.Range("A174").Value = "59082104"
$ws.Range("B174").Value = "Player-59082104"
$ws.Range("C174").Value = ""
$ws.Range("D174").Value = ""
$ws.Range("E174").Value = "一馆"
$ws.Range("CT174").Value = 25
$ws.Range("CT174").Interior.Color = 16777215
$ws.Range("CU174").NumberFormat = "@"
$ws.Range("CU174").Value = "3994"
